# Update gh-pages to output generated at 456a3b4
# Bumps the "想去人数" (F column) counts on several rows across the
# "展览", "演出", "本地生活" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value  = 7975
$ws.Range("F10").Value = 1960
$ws.Range("F12").Value = 195
$ws.Range("F13").Value = 1834
$ws.Range("F25").Value = 1224
$ws.Range("F36").Value = 3660
$ws.Range("F43").Value = 753

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 84
$ws.Range("F20").Value = 60
$ws.Range("F32").Value = 64

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 1429
$ws.Range("F9").Value = 9105

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value  = 7975
$ws.Range("F10").Value = 195
$ws.Range("F11").Value = 1834
$ws.Range("F23").Value = 1224
$ws.Range("F34").Value = 84
$ws.Range("F36").Value = 3660
$ws.Range("F42").Value = 753
